# Auto-generated edit script: updates Coin/Link/Price/Volume(1h) cells on Sheet1
# per the "Updated symbol list" commit (crypto price-ticker refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: cell reference, new text value.
# Numeric-looking values (prices / percentages) are given a leading
# apostrophe so Excel stores them as literal text (matching the workbook's
# existing convention of text-formatted Price/Volume columns) instead of
# silently converting them to numbers.
$cellUpdates = @(
    @{ Cell = "D2"; Value = '''297.28' }
    @{ Cell = "E2"; Value = '''1.66%' }
    @{ Cell = "D3"; Value = '''41.85' }
    @{ Cell = "E3"; Value = '''3.58%' }
    @{ Cell = "D4"; Value = '''5.019' }
    @{ Cell = "E4"; Value = '''-0.37%' }
    @{ Cell = "D5"; Value = '''0.07526' }
    @{ Cell = "E5"; Value = '''2.76%' }
    @{ Cell = "B6"; Value = 'GateToken' }
    @{ Cell = "C6"; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' }
    @{ Cell = "D6"; Value = '''4.371' }
    @{ Cell = "E6"; Value = '''1.86%' }
    @{ Cell = "B7"; Value = 'FTXToken' }
    @{ Cell = "C7"; Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt' }
    @{ Cell = "D7"; Value = '''1.583' }
    @{ Cell = "E7"; Value = '''3.69%' }
    @{ Cell = "B8"; Value = 'MXToken' }
    @{ Cell = "C8"; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' }
    @{ Cell = "D8"; Value = '''0.9231' }
    @{ Cell = "E8"; Value = '''-0.54%' }
    @{ Cell = "B9"; Value = 'BTSEToken' }
    @{ Cell = "C9"; Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse' }
    @{ Cell = "D9"; Value = '''2.401' }
    @{ Cell = "E9"; Value = '''2.00%' }
    @{ Cell = "B10"; Value = 'LiechtensteinCryptoassetsExchange' }
    @{ Cell = "C10"; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' }
    @{ Cell = "D10"; Value = '''0.1192' }
    @{ Cell = "E10"; Value = '''1.34%' }
    @{ Cell = "B11"; Value = 'WazirX' }
    @{ Cell = "C11"; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' }
    @{ Cell = "D11"; Value = '''0.1829' }
    @{ Cell = "E11"; Value = '''4.72%' }
    @{ Cell = "B12"; Value = 'MandalaExchangeToken' }
    @{ Cell = "C12"; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' }
    @{ Cell = "D12"; Value = '''0.08951' }
    @{ Cell = "E12"; Value = '''3.58%' }
    @{ Cell = "B13"; Value = 'BitrueCoin' }
    @{ Cell = "C13"; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' }
    @{ Cell = "D13"; Value = '''0.04092' }
    @{ Cell = "E13"; Value = '''-5.50%' }
    @{ Cell = "B14"; Value = 'BitMartToken' }
    @{ Cell = "C14"; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' }
    @{ Cell = "D14"; Value = '''0.1050' }
    @{ Cell = "E14"; Value = '''-0.54%' }
    @{ Cell = "B15"; Value = 'BitForexToken' }
    @{ Cell = "C15"; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' }
    @{ Cell = "D15"; Value = '''0.001277' }
    @{ Cell = "E15"; Value = '''0.27%' }
    @{ Cell = "B16"; Value = 'TigerCash' }
    @{ Cell = "C16"; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' }
    @{ Cell = "D16"; Value = '''0.005872' }
    @{ Cell = "E16"; Value = '''1.74%' }
    @{ Cell = "B17"; Value = 'LEO' }
    @{ Cell = "C17"; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' }
    @{ Cell = "D17"; Value = '''3.342' }
    @{ Cell = "E17"; Value = '''0.14%' }
    @{ Cell = "D18"; Value = '''0.3314' }
    @{ Cell = "E18"; Value = '''0.79%' }
    @{ Cell = "D19"; Value = '''8.086' }
    @{ Cell = "E19"; Value = '''1.48%' }
    @{ Cell = "D20"; Value = '''0.1381' }
    @{ Cell = "E20"; Value = '''-0.68%' }
    @{ Cell = "E21"; Value = '''13.14%' }
    @{ Cell = "D22"; Value = '''0.04085' }
    @{ Cell = "E22"; Value = '''3.76%' }
    @{ Cell = "E23"; Value = '''0.40%' }
    @{ Cell = "D24"; Value = '''0.003896' }
    @{ Cell = "E24"; Value = '''2.93%' }
    @{ Cell = "E25"; Value = '''-3.90%' }
    @{ Cell = "D38"; Value = '''0.02413' }
    @{ Cell = "E38"; Value = '''5.58%' }
    @{ Cell = "D39"; Value = '''0.05214' }
    @{ Cell = "E39"; Value = '''3.47%' }
    @{ Cell = "D40"; Value = '''0.006304' }
    @{ Cell = "E40"; Value = '''5.72%' }
    @{ Cell = "D41"; Value = '''0.007810' }
    @{ Cell = "E41"; Value = '''1.55%' }
    @{ Cell = "D42"; Value = '''0.1328' }
    @{ Cell = "E42"; Value = '''3.20%' }
    @{ Cell = "D43"; Value = '''0.007392' }
    @{ Cell = "E43"; Value = '''0.59%' }
    @{ Cell = "D44"; Value = '''0.007774' }
    @{ Cell = "E44"; Value = '''-6.02%' }
    @{ Cell = "D45"; Value = '''0.2966' }
    @{ Cell = "E45"; Value = '''1.73%' }
    @{ Cell = "D46"; Value = '''0.00006588' }
    @{ Cell = "E46"; Value = '''5.15%' }
    @{ Cell = "D47"; Value = '''0.00000000750' }
    @{ Cell = "E47"; Value = '''-0.01%' }
    @{ Cell = "D48"; Value = '''0.03172' }
    @{ Cell = "E48"; Value = '''0.93%' }
    @{ Cell = "D49"; Value = '''0.004203' }
    @{ Cell = "E49"; Value = '''0.02%' }
    @{ Cell = "D50"; Value = '''0.00002101' }
    @{ Cell = "E50"; Value = '''-0.01%' }
    @{ Cell = "D51"; Value = '''0.0002001' }
    @{ Cell = "E51"; Value = '''-0.01%' }
)

foreach ($update in $cellUpdates) {
    $ws.Range($update.Cell).Value = $update.Value
}
